$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "318.19"
Set-TextValue $ws.Range("E2") "3.85%"

# Row 3
Set-TextValue $ws.Range("D3") "39.88"
Set-TextValue $ws.Range("E3") "2.04%"

# Row 4
Set-TextValue $ws.Range("D4") "5.140"
Set-TextValue $ws.Range("E4") "0.89%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08214"

# Row 6
Set-TextValue $ws.Range("D6") "2.050"
Set-TextValue $ws.Range("E6") "4.99%"

# Row 7
Set-TextValue $ws.Range("D7") "8.318"
Set-TextValue $ws.Range("E7") "3.96%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9342"
Set-TextValue $ws.Range("E8") "0.21%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D9") "0.1361"
Set-TextValue $ws.Range("E9") "-6.06%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1981"
Set-TextValue $ws.Range("E10") "2.66%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.09058"
Set-TextValue $ws.Range("E11") "0.72%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03475"
Set-TextValue $ws.Range("E12") "-0.95%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09800"
Set-TextValue $ws.Range("E13") "0.08%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001406"
Set-TextValue $ws.Range("E14") "0.34%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.006313"
Set-TextValue $ws.Range("E15") "4.64%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.681"
Set-TextValue $ws.Range("E16") "-2.88%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "4.294"
Set-TextValue $ws.Range("E17") "2.38%"

# Row 18
Set-TextValue $ws.Range("E18") "-2.67%"

# Row 19
Set-TextValue $ws.Range("D19") "0.3475"
Set-TextValue $ws.Range("E19") "1.51%"

# Row 20
Set-TextValue $ws.Range("D20") "0.1297"
Set-TextValue $ws.Range("E20") "-0.52%"

# Row 21
Set-TextValue $ws.Range("D21") "4.889"
Set-TextValue $ws.Range("E21") "7.26%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2445"
Set-TextValue $ws.Range("E22") "1.17%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04321"
Set-TextValue $ws.Range("E23") "-1.25%"

# Row 24
Set-TextValue $ws.Range("D24") "0.001225"
Set-TextValue $ws.Range("E24") "-1.05%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004772"
Set-TextValue $ws.Range("E25") "11.72%"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001297"
Set-TextValue $ws.Range("E26") "-0.43%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0003990"
Set-TextValue $ws.Range("E27") "-10.29%"

# Row 39
Set-TextValue $ws.Range("D39") "0.02225"
Set-TextValue $ws.Range("E39") "9.87%"

# Row 40
Set-TextValue $ws.Range("D40") "0.05221"
Set-TextValue $ws.Range("E40") "3.06%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007607"
Set-TextValue $ws.Range("E41") "2.11%"

# Row 42
Set-TextValue $ws.Range("D42") "0.009743"
Set-TextValue $ws.Range("E42") "-5.33%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1387"
Set-TextValue $ws.Range("E43") "2.94%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002088"
Set-TextValue $ws.Range("E44") "-1.67%"

# Row 45
Set-TextValue $ws.Range("D45") "0.009181"
Set-TextValue $ws.Range("E45") "0.70%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006551"
Set-TextValue $ws.Range("E46") "5.73%"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000748"
Set-TextValue $ws.Range("E47") "-0.48%"

# Row 48
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D48") "0.002982"
Set-TextValue $ws.Range("E48") "-3.68%"

# Row 49
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D49") "0.001686"
Set-TextValue $ws.Range("E49") "5.25%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002095"
Set-TextValue $ws.Range("E50") "-0.48%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0001995"
Set-TextValue $ws.Range("E51") "-0.48%"
